# Updates cryptos list: price (D) and volume/1h change (E) columns,
# plus a full row replacement (row 51: Algorand -> EnergySwap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '27.513.21'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.37%  '

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.645.50'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -1.21%  '

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '212.41'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.50%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.531'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +4.14%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '23.50'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -2.30%  '

$ws.Range('E9').Value = '  -2.08%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0611'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -1.41%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0892'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.45%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.878.88'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -1.18%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.640.41'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -1.29%  '

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.590'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +3.37%  '

$ws.Range('E15').Value = '  -2.15%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '64.52'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -2.79%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '27.478.36'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.44%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '231.29'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -4.34%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.0₃0724'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.95%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.56'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.36%  '

$ws.Range('E21').Value = '  -0.07%  '

$ws.Range('E22').Value = '  -3.83%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '9.74'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +4.08%  '

$ws.Range('E24').Value = '  -1.00%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '148.11'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.94%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '7.04'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -3.20%  '

$ws.Range('E27').Value = '  +1.81%  '

$ws.Range('E28').Value = '  -0.05%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '15.66'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -4.31%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.18'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -3.83%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.32'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.98%  '

$ws.Range('E33').Value = '  +1.54%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.424.76'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -2.52%  '

$ws.Range('E35').Value = '  +0.30%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.38'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.30%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.570'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.28%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.890'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -4.47%  '

$ws.Range('E39').Value = '  -3.30%  '

$ws.Range('E40').Value = '  -1.02%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.05%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.819'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.36%  '

$ws.Range('E43').Value = '  +2.59%  '

$ws.Range('E44').Value = '  -1.79%  '

$ws.Range('E45').Value = '  +1.01%  '

$ws.Range('E46').Value = '  -7.16%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.788.44'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.13%  '

$ws.Range('E48').Value = '  -2.83%  '

$ws.Range('E49').Value = '  -0.45%  '

$ws.Range('E50').Value = '  -0.93%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '7.79'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.60%  '
